$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 ("1.00") is never modified by this edit and carries the default (unstyled) cell
# format; we reuse its .Style below to strip the quote-prefix formatting that Excel
# applies automatically when a numeric-looking string is assigned, so number-like
# price text (e.g. "508.11") is stored as text without altering cell style.
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "57.484.24"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.409.52"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'508.11"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").Value = "'133.07"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'0.558"
$ws.Range("D9").Value = "2.446.43"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'0.0980"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  -7.83%  "
$ws.Range("D14").Value = "2.844.12"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "57.327.66"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "'21.87"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "2.402.08"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").Value = "'10.29"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'313.94"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'6.42"
$ws.Range("E22").Value = "  +5.19%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").Value = "'65.13"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "2.523.38"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").Value = "'7.55"
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("D31").Value = "'173.69"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "0.0₃0733"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "'6.17"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'18.00"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("D40").Value = "'3.83"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'0.815"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "'135.04"
$ws.Range("E44").Value = "  +10.33%  "
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D47").Value = "'255.98"
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("D48").Value = "'0.572"
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").Value = "'0.0917"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "'0.0493"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  +0.46%  "

# Reset style on cells that received quote-prefixed numeric-looking text so their
# cell formatting matches the rest of the sheet (no lingering quote-prefix style).
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").Style = $plainStyle
$ws.Range("D7").Style = $plainStyle
$ws.Range("D8").Style = $plainStyle
$ws.Range("D10").Style = $plainStyle
$ws.Range("D16").Style = $plainStyle
$ws.Range("D19").Style = $plainStyle
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").Style = $plainStyle
$ws.Range("D23").Style = $plainStyle
$ws.Range("D25").Style = $plainStyle
$ws.Range("D30").Style = $plainStyle
$ws.Range("D31").Style = $plainStyle
$ws.Range("D34").Style = $plainStyle
$ws.Range("D38").Style = $plainStyle
$ws.Range("D40").Style = $plainStyle
$ws.Range("D42").Style = $plainStyle
$ws.Range("D44").Style = $plainStyle
$ws.Range("D47").Style = $plainStyle
$ws.Range("D48").Style = $plainStyle
$ws.Range("D49").Style = $plainStyle
$ws.Range("D50").Style = $plainStyle
